# "Actualización desde MV -datos-"
# Appends the new quarterly data row (01-04-2021) to the bottom of the
# "Balance Contable BCCh" table on Sheet1 (row 63, right after the existing
# last row, 62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63

# Column A ("Serie") holds quarter-start dates stored as plain text (every
# other cell in the column is a shared string like "01-01-2021", not a real
# Excel date serial). Force the cell to text first so the "01-04-2021"
# value isn't auto-converted into a date serial, then restore the cell's
# style to match the rest of the column (plain/default formatting, same as
# row 62) so no stray per-cell formatting is left behind.
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("A" + $row).Value = "01-04-2021"
$ws.Range("A" + $row).Style = $ws.Range("A62").Style

$ws.Range("B" + $row).Value = 97712
$ws.Range("C" + $row).Value = 0
$ws.Range("D" + $row).Value = 44954
$ws.Range("E" + $row).Value = 0
$ws.Range("F" + $row).Value = 51877
$ws.Range("G" + $row).Value = 880
$ws.Range("H" + $row).Value = 99378
$ws.Range("I" + $row).Value = 43335
$ws.Range("J" + $row).Value = 42794
$ws.Range("K" + $row).Value = 541
$ws.Range("L" + $row).Value = 0
$ws.Range("M" + $row).Value = 0
$ws.Range("N" + $row).Value = 7747
$ws.Range("O" + $row).Value = 43786
$ws.Range("P" + $row).Value = 4510
$ws.Range("Q" + $row).Value = -1666
$ws.Range("R" + $row).Value = 1329
$ws.Range("S" + $row).Value = 40010
